$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 5.488907176552729

$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 14.40014219143469

$ws.Range("B4").Value = 0.0003714022599530242
$ws.Range("C4").Value = 0.0001537489499301437
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 0.6533176403021212

$ws.Range("B5").Value = 0.3464964993005633
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 16.98373111632243
$ws.Range("E5").Value = 71517.89157740913
$ws.Range("G5").Value = 71536.87504148365

$ws.Range("B6").Value = 3.182878228561681
$ws.Range("C6").Value = 9.226618575922256
$ws.Range("D6").Value = 0.7127328510149897
$ws.Range("E6").Value = 6.48142807727062
$ws.Range("G6").Value = 19.60365773276954

$ws.Range("B7").Value = 0.001754667048134761
$ws.Range("C7").Value = 0.3375848360084654
$ws.Range("D7").Value = 0.7127328510149897
$ws.Range("E7").Value = 6.48142807727062
$ws.Range("G7").Value = 7.53350043134221
